$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.4483985765124555
$ws1.Range("C2").Value = 0.07530120481927711
$ws1.Range("D2").Value = 0.8928571428571429
$ws1.Range("E2").Value = 0.1388888888888889
$ws1.Range("F2").Value = 0.2815315315315315
$ws1.Range("G2").Value = 0.6298449612403101
$ws1.Range("H2").Value = 0.7532771535580525
$ws1.Range("I2").Value = 25
$ws1.Range("J2").Value = 307
$ws1.Range("K2").Value = 227
$ws1.Range("L2").Value = 3

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9869565217391304
$ws2.Range("C2").Value = 0.4250936329588015
$ws2.Range("D2").Value = 0.5942408376963351

$ws2.Range("B3").Value = 0.07530120481927711
$ws2.Range("C3").Value = 0.8928571428571429
$ws2.Range("D3").Value = 0.1388888888888889

$ws2.Range("B4").Value = 0.4483985765124555
$ws2.Range("C4").Value = 0.4483985765124555
$ws2.Range("D4").Value = 0.4483985765124555
$ws2.Range("E4").Value = 0.4483985765124555

$ws2.Range("B5").Value = 0.5311288632792037
$ws2.Range("C5").Value = 0.6589753879079722
$ws2.Range("D5").Value = 0.366564863292612

$ws2.Range("B6").Value = 0.9415359721417002
$ws2.Range("C6").Value = 0.4483985765124555
$ws2.Range("D6").Value = 0.5715542637343983

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 227
$ws3.Range("C2").Value = 307
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 25
